$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: B1 "name" -> "r"
$ws.Range("B1").Value = "r"

# Update row 2: B2 "asd" -> "ew"
$ws.Range("B2").Value = "ew"

# Add new row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "s"
$ws.Range("C4").Value = "s"

# Update active selection to reflect the last edited cell (C4)
$ws.Range("C4").Select() | Out-Null
